$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -12.965
$ws.Range("B9").Value = 5.317
$ws.Range("C12").Value = -11.517
$ws.Range("E13").Value = 16.532
$ws.Range("D15").Value = -8.463000000000001
$ws.Range("E16").Value = 16.643
$ws.Range("B18").Value = 5.282999999999999
$ws.Range("B20").Value = 7.592000000000001
$ws.Range("E20").Value = 16.076
$ws.Range("E24").Value = 16.83
$ws.Range("C26").Value = -12.808
$ws.Range("B27").Value = 5.583
$ws.Range("C27").Value = -13.501
$ws.Range("C29").Value = -12.335
$ws.Range("C37").Value = -13.351
$ws.Range("C38").Value = -13.738
$ws.Range("D38").Value = -7.473000000000001
$ws.Range("E39").Value = 16.373
$ws.Range("D44").Value = -7.568999999999998
$ws.Range("E48").Value = 17.347
$ws.Range("C51").Value = -12.613
$ws.Range("D51").Value = -7.626
$ws.Range("E52").Value = 16.826
$ws.Range("C55").Value = -13.65
$ws.Range("E56").Value = 16.748
$ws.Range("D57").Value = -8.032
$ws.Range("D63").Value = -7.337000000000001
$ws.Range("B69").Value = 5.667
$ws.Range("C69").Value = -11.17
$ws.Range("C70").Value = -13.202
$ws.Range("D70").Value = -7.992999999999999
$ws.Range("B76").Value = 6.723999999999999
$ws.Range("B82").Value = 5.345000000000001
$ws.Range("C83").Value = -13.509
$ws.Range("E84").Value = 16.68
$ws.Range("D99").Value = -7.512
$ws.Range("E100").Value = 16.569
$ws.Range("E101").Value = 16.821
$ws.Range("C102").Value = -13.419
